# edit.ps1 -- apply the changes described by the diff:
#   1. "Version 0.8.0" -> "Version 0.8.2"
#   2. Insert a comment ("This is a comment.") anchored on the text
#      "(3) First sentence in a new paragraph. " (the italic trailing
#      space run included) in the third numbered paragraph.
#   3. "an import work" -> "an important book" (inside the larger
#      sentence "And here is a reference to an import work:" ->
#      "And here is a reference to an important book:")

$d = $word.ActiveDocument

# --- 1. Version bump -------------------------------------------------
$d.Content.Find.Execute("Version 0.8.0", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Version 0.8.2", 2) | Out-Null

# --- 2. Add the comment ------------------------------------------------
# Locate "(3)" that starts the paragraph we need to annotate, and the
# "Second link." hyperlink text that immediately follows the comment
# range in that same paragraph, so we can build a precise Range that
# covers "(3) First sentence in a new paragraph. " (trailing italic
# space included), matching commentRangeStart/commentRangeEnd placement.

$startMarker = $d.Content.Duplicate
$startMarker.Find.Execute("(3)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null

$endMarker = $d.Content.Duplicate
$endMarker.Find.Execute("Second link.", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null

$commentRange = $d.Range($startMarker.Start, $endMarker.Start)

$word.UserName = "Per Kraulis"
$word.UserInitials = ""

$comment = $d.Comments.Add($commentRange, "This is a comment.")

# --- 3. Fix the reference to Darwin's book -----------------------------
$d.Content.Find.Execute("an import work", $true, $false, $false, $false, `
    $false, $true, 1, $false, "an important book", 2) | Out-Null
